$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New column BJ: "quit" header + description, matching the existing
# header-row / description-row styling (wrap text, same row heights).
$ws.Range("BJ1").Value = "quit"
$ws.Range("BJ2").Value = "sub quit the experiment before reaching this trial"
$ws.Range("BJ1:BJ2").WrapText = $true

# Row heights shrink slightly as part of the same formatting pass.
$ws.Rows.Item(1).RowHeight = 27.6
$ws.Rows.Item(2).RowHeight = 69

# Move the selection/view to the newly added last column, like the
# author scrolling over to see their new column after adding it.
$ws.Range("BJ3").Select()
